$d = $word.ActiveDocument
$d.Content.Find.Execute(" RD", $true, $false, $false, $false, $false,
                         $true, 1, $false, " ARD", 2)
